# ============================================================
# 1) Restructure worksheets:
#    before : [ODI Batting]
#    after  : [Player Info, ODI Batting, ODI Batting Extra]
#
# NOTE: worksheet handles in this host resolve by POSITIONAL
# index at the time they're dereferenced, not by stable identity.
# Inserting a sheet BEFORE a held index shifts what that index
# refers to. So: do the 'insert after' first (doesn't disturb the
# original sheet's index), then the 'insert before' using a freshly
# fetched reference, and finally re-fetch ALL sheets by their
# settled final index before writing any cell data.
# ============================================================
$wb = $excel.ActiveWorkbook

$battingOrig = $wb.Worksheets.Item(1)

# Insert the new last sheet right after 'ODI Batting' (index unaffected)
$extraNew = $wb.Worksheets.Add($null, $battingOrig)
$extraNew.Name = "ODI Batting Extra"

# Insert the new first sheet right before 'ODI Batting' (re-fetch fresh first)
$battingForInsert = $wb.Worksheets.Item(1)
$infoNew = $wb.Worksheets.Add($battingForInsert)
$infoNew.Name = "Player Info"

# Settle final, stable references by final position:
#   1 = Player Info, 2 = ODI Batting, 3 = ODI Batting Extra
$playerInfo  = $wb.Worksheets.Item(1)
$batting     = $wb.Worksheets.Item(2)
$battingExtra = $wb.Worksheets.Item(3)

# ============================================================
# 2) Populate 'Player Info'
# ============================================================
$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

$playerInfo.Cells.Item(2,1).NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = "4542"
$playerInfo.Cells.Item(2,2).Value = "Rishabh Rajendra Pant"
$playerInfo.Cells.Item(2,3).Value = "Left Handed"
$playerInfo.Cells.Item(2,4).Value = "Does Not Bowl | Unknown"

# ============================================================
# 3) Update 'ODI Batting':
#    - header D1: MATCH_CARD_LINK -> MATCH_CODE
#    - column D: full scorecard URL -> bare numeric match code
#    - drop the stray empty INNING_NUMBER cells (col B) on rows
#      2, 11, 26 and 30 (these were 'did not bat' rows)
# ============================================================
$batting.Cells.Item(1,4).Value = "MATCH_CODE"

$batting.Cells.Item(2,4).NumberFormat = "@"
$batting.Cells.Item(2,4).Value = "4213"
$batting.Cells.Item(3,4).NumberFormat = "@"
$batting.Cells.Item(3,4).Value = "4216"
$batting.Cells.Item(4,4).NumberFormat = "@"
$batting.Cells.Item(4,4).Value = "4219"
$batting.Cells.Item(5,4).NumberFormat = "@"
$batting.Cells.Item(5,4).Value = "4268"
$batting.Cells.Item(6,4).NumberFormat = "@"
$batting.Cells.Item(6,4).Value = "4270"
$batting.Cells.Item(7,4).NumberFormat = "@"
$batting.Cells.Item(7,4).Value = "4342"
$batting.Cells.Item(8,4).NumberFormat = "@"
$batting.Cells.Item(8,4).Value = "4345"
$batting.Cells.Item(9,4).NumberFormat = "@"
$batting.Cells.Item(9,4).Value = "4350"
$batting.Cells.Item(10,4).NumberFormat = "@"
$batting.Cells.Item(10,4).Value = "4353"
$batting.Cells.Item(11,4).NumberFormat = "@"
$batting.Cells.Item(11,4).Value = "4359"
$batting.Cells.Item(12,4).NumberFormat = "@"
$batting.Cells.Item(12,4).Value = "4360"
$batting.Cells.Item(13,4).NumberFormat = "@"
$batting.Cells.Item(13,4).Value = "4362"
$batting.Cells.Item(14,4).NumberFormat = "@"
$batting.Cells.Item(14,4).Value = "4385"
$batting.Cells.Item(15,4).NumberFormat = "@"
$batting.Cells.Item(15,4).Value = "4387"
$batting.Cells.Item(16,4).NumberFormat = "@"
$batting.Cells.Item(16,4).Value = "4388"
$batting.Cells.Item(17,4).NumberFormat = "@"
$batting.Cells.Item(17,4).Value = "4398"
$batting.Cells.Item(18,4).NumberFormat = "@"
$batting.Cells.Item(18,4).Value = "4456"
$batting.Cells.Item(19,4).NumberFormat = "@"
$batting.Cells.Item(19,4).Value = "4457"
$batting.Cells.Item(20,4).NumberFormat = "@"
$batting.Cells.Item(20,4).Value = "4524"
$batting.Cells.Item(21,4).NumberFormat = "@"
$batting.Cells.Item(21,4).Value = "4526"
$batting.Cells.Item(22,4).NumberFormat = "@"
$batting.Cells.Item(22,4).Value = "4529"
$batting.Cells.Item(23,4).NumberFormat = "@"
$batting.Cells.Item(23,4).Value = "4533"
$batting.Cells.Item(24,4).NumberFormat = "@"
$batting.Cells.Item(24,4).Value = "4535"
$batting.Cells.Item(25,4).NumberFormat = "@"
$batting.Cells.Item(25,4).Value = "4536"
$batting.Cells.Item(26,4).NumberFormat = "@"
$batting.Cells.Item(26,4).Value = "4609"
$batting.Cells.Item(27,4).NumberFormat = "@"
$batting.Cells.Item(27,4).Value = "4613"
$batting.Cells.Item(28,4).NumberFormat = "@"
$batting.Cells.Item(28,4).Value = "4618"
$batting.Cells.Item(29,4).NumberFormat = "@"
$batting.Cells.Item(29,4).Value = "4669"
$batting.Cells.Item(30,4).NumberFormat = "@"
$batting.Cells.Item(30,4).Value = "4673"
$batting.Cells.Item(31,4).NumberFormat = "@"
$batting.Cells.Item(31,4).Value = "4676"

$batting.Cells.Item(2,2).ClearContents()
$batting.Cells.Item(11,2).ClearContents()
$batting.Cells.Item(26,2).ClearContents()
$batting.Cells.Item(30,2).ClearContents()

# ============================================================
# 4) Populate 'ODI Batting Extra'
# ============================================================
$battingExtra.Cells.Item(1,1).Value = "MATCH_CODE"
$battingExtra.Cells.Item(1,2).Value = "BATTING_POSITION"
$battingExtra.Cells.Item(1,3).Value = "NUM_4"
$battingExtra.Cells.Item(1,4).Value = "NUM_6"
$battingExtra.Cells.Item(1,5).Value = "PERCENT_RUNS_OF_TOTAL"
$battingExtra.Cells.Item(1,6).Value = "MAN_OF_MATCH"

# Row 2: match 4360
$battingExtra.Cells.Item(2,1).NumberFormat = "@"
$battingExtra.Cells.Item(2,1).Value = "4360"
$battingExtra.Cells.Item(2,2).Value = 4
$battingExtra.Cells.Item(2,3).NumberFormat = "@"
$battingExtra.Cells.Item(2,3).Value = "2"
$battingExtra.Cells.Item(2,4).NumberFormat = "@"
$battingExtra.Cells.Item(2,4).Value = "0"
$battingExtra.Cells.Item(2,5).Value = "7.17%"
$battingExtra.Cells.Item(2,6).Value = "NO"

# Row 3: match 4362
$battingExtra.Cells.Item(3,1).NumberFormat = "@"
$battingExtra.Cells.Item(3,1).Value = "4362"
$battingExtra.Cells.Item(3,2).Value = ""
$battingExtra.Cells.Item(3,3).Value = ""
$battingExtra.Cells.Item(3,4).Value = ""
$battingExtra.Cells.Item(3,5).Value = ""
$battingExtra.Cells.Item(3,6).Value = "NO"

# Row 4: match 4385
$battingExtra.Cells.Item(4,1).NumberFormat = "@"
$battingExtra.Cells.Item(4,1).Value = "4385"
$battingExtra.Cells.Item(4,2).Value = ""
$battingExtra.Cells.Item(4,3).Value = ""
$battingExtra.Cells.Item(4,4).Value = ""
$battingExtra.Cells.Item(4,5).Value = ""
$battingExtra.Cells.Item(4,6).Value = "NO"

# Row 5: match 4387
$battingExtra.Cells.Item(5,1).NumberFormat = "@"
$battingExtra.Cells.Item(5,1).Value = "4387"
$battingExtra.Cells.Item(5,2).Value = 5
$battingExtra.Cells.Item(5,3).NumberFormat = "@"
$battingExtra.Cells.Item(5,3).Value = "3"
$battingExtra.Cells.Item(5,4).NumberFormat = "@"
$battingExtra.Cells.Item(5,4).Value = "4"
$battingExtra.Cells.Item(5,5).Value = "10.08%"
$battingExtra.Cells.Item(5,6).Value = "NO"

# Row 6: match 4388
$battingExtra.Cells.Item(6,1).NumberFormat = "@"
$battingExtra.Cells.Item(6,1).Value = "4388"
$battingExtra.Cells.Item(6,2).Value = 5
$battingExtra.Cells.Item(6,3).NumberFormat = "@"
$battingExtra.Cells.Item(6,3).Value = "1"
$battingExtra.Cells.Item(6,4).NumberFormat = "@"
$battingExtra.Cells.Item(6,4).Value = "0"
$battingExtra.Cells.Item(6,5).Value = "2.22%"
$battingExtra.Cells.Item(6,6).Value = "NO"

# Row 7: match 4398
$battingExtra.Cells.Item(7,1).NumberFormat = "@"
$battingExtra.Cells.Item(7,1).Value = "4398"
$battingExtra.Cells.Item(7,2).Value = 6
$battingExtra.Cells.Item(7,3).NumberFormat = "@"
$battingExtra.Cells.Item(7,3).Value = "2"
$battingExtra.Cells.Item(7,4).NumberFormat = "@"
$battingExtra.Cells.Item(7,4).Value = "1"
$battingExtra.Cells.Item(7,5).Value = "10.98%"
$battingExtra.Cells.Item(7,6).Value = "NO"

# Row 8: match 4456
$battingExtra.Cells.Item(8,1).NumberFormat = "@"
$battingExtra.Cells.Item(8,1).Value = "4456"
$battingExtra.Cells.Item(8,2).Value = ""
$battingExtra.Cells.Item(8,3).Value = ""
$battingExtra.Cells.Item(8,4).Value = ""
$battingExtra.Cells.Item(8,5).Value = ""
$battingExtra.Cells.Item(8,6).Value = "NO"

# Row 9: match 4457
$battingExtra.Cells.Item(9,1).NumberFormat = "@"
$battingExtra.Cells.Item(9,1).Value = "4457"
$battingExtra.Cells.Item(9,2).Value = 4
$battingExtra.Cells.Item(9,3).NumberFormat = "@"
$battingExtra.Cells.Item(9,3).Value = "5"
$battingExtra.Cells.Item(9,4).NumberFormat = "@"
$battingExtra.Cells.Item(9,4).Value = "4"
$battingExtra.Cells.Item(9,5).Value = "23.71%"
$battingExtra.Cells.Item(9,6).Value = "NO"

# Row 10: match 4524
$battingExtra.Cells.Item(10,1).NumberFormat = "@"
$battingExtra.Cells.Item(10,1).Value = "4524"
$battingExtra.Cells.Item(10,2).Value = ""
$battingExtra.Cells.Item(10,3).Value = ""
$battingExtra.Cells.Item(10,4).Value = ""
$battingExtra.Cells.Item(10,5).Value = ""
$battingExtra.Cells.Item(10,6).Value = "NO"

# Row 11: match 4526
$battingExtra.Cells.Item(11,1).NumberFormat = "@"
$battingExtra.Cells.Item(11,1).Value = "4526"
$battingExtra.Cells.Item(11,2).Value = 4
$battingExtra.Cells.Item(11,3).NumberFormat = "@"
$battingExtra.Cells.Item(11,3).Value = "10"
$battingExtra.Cells.Item(11,4).NumberFormat = "@"
$battingExtra.Cells.Item(11,4).Value = "2"
$battingExtra.Cells.Item(11,5).Value = "29.62%"
$battingExtra.Cells.Item(11,6).Value = "NO"

# Row 12: match 4529
$battingExtra.Cells.Item(12,1).NumberFormat = "@"
$battingExtra.Cells.Item(12,1).Value = "4529"
$battingExtra.Cells.Item(12,2).Value = 4
$battingExtra.Cells.Item(12,3).NumberFormat = "@"
$battingExtra.Cells.Item(12,3).Value = "0"
$battingExtra.Cells.Item(12,4).NumberFormat = "@"
$battingExtra.Cells.Item(12,4).Value = "0"
$battingExtra.Cells.Item(12,5).Value = ""
$battingExtra.Cells.Item(12,6).Value = "NO"

# Row 13: match 4533
$battingExtra.Cells.Item(13,1).NumberFormat = "@"
$battingExtra.Cells.Item(13,1).Value = "4533"
$battingExtra.Cells.Item(13,2).Value = ""
$battingExtra.Cells.Item(13,3).Value = ""
$battingExtra.Cells.Item(13,4).Value = ""
$battingExtra.Cells.Item(13,5).Value = ""
$battingExtra.Cells.Item(13,6).Value = "NO"

# Row 14: match 4535
$battingExtra.Cells.Item(14,1).NumberFormat = "@"
$battingExtra.Cells.Item(14,1).Value = "4535"
$battingExtra.Cells.Item(14,2).Value = 2
$battingExtra.Cells.Item(14,3).NumberFormat = "@"
$battingExtra.Cells.Item(14,3).Value = "3"
$battingExtra.Cells.Item(14,4).NumberFormat = "@"
$battingExtra.Cells.Item(14,4).Value = "0"
$battingExtra.Cells.Item(14,5).Value = "7.59%"
$battingExtra.Cells.Item(14,6).Value = "NO"

# Row 15: match 4536
$battingExtra.Cells.Item(15,1).NumberFormat = "@"
$battingExtra.Cells.Item(15,1).Value = "4536"
$battingExtra.Cells.Item(15,2).Value = 5
$battingExtra.Cells.Item(15,3).NumberFormat = "@"
$battingExtra.Cells.Item(15,3).Value = "6"
$battingExtra.Cells.Item(15,4).NumberFormat = "@"
$battingExtra.Cells.Item(15,4).Value = "1"
$battingExtra.Cells.Item(15,5).Value = "21.13%"
$battingExtra.Cells.Item(15,6).Value = "NO"

# Row 16: match 4609
$battingExtra.Cells.Item(16,1).NumberFormat = "@"
$battingExtra.Cells.Item(16,1).Value = "4609"
$battingExtra.Cells.Item(16,2).Value = 5
$battingExtra.Cells.Item(16,3).Value = ""
$battingExtra.Cells.Item(16,4).Value = ""
$battingExtra.Cells.Item(16,5).Value = ""
$battingExtra.Cells.Item(16,6).Value = "NO"

# Row 17: match 4613
$battingExtra.Cells.Item(17,1).NumberFormat = "@"
$battingExtra.Cells.Item(17,1).Value = "4613"
$battingExtra.Cells.Item(17,2).Value = 4
$battingExtra.Cells.Item(17,3).NumberFormat = "@"
$battingExtra.Cells.Item(17,3).Value = "0"
$battingExtra.Cells.Item(17,4).NumberFormat = "@"
$battingExtra.Cells.Item(17,4).Value = "0"
$battingExtra.Cells.Item(17,5).Value = ""
$battingExtra.Cells.Item(17,6).Value = "NO"

# Row 18: match 4618
$battingExtra.Cells.Item(18,1).NumberFormat = "@"
$battingExtra.Cells.Item(18,1).Value = "4618"
$battingExtra.Cells.Item(18,2).Value = 4
$battingExtra.Cells.Item(18,3).NumberFormat = "@"
$battingExtra.Cells.Item(18,3).Value = "16"
$battingExtra.Cells.Item(18,4).NumberFormat = "@"
$battingExtra.Cells.Item(18,4).Value = "2"
$battingExtra.Cells.Item(18,5).Value = "47.89%"
$battingExtra.Cells.Item(18,6).Value = "YES"

# Row 19: match 4669
$battingExtra.Cells.Item(19,1).NumberFormat = "@"
$battingExtra.Cells.Item(19,1).Value = "4669"
$battingExtra.Cells.Item(19,2).Value = 4
$battingExtra.Cells.Item(19,3).NumberFormat = "@"
$battingExtra.Cells.Item(19,3).Value = "2"
$battingExtra.Cells.Item(19,4).NumberFormat = "@"
$battingExtra.Cells.Item(19,4).Value = "0"
$battingExtra.Cells.Item(19,5).Value = "4.90%"
$battingExtra.Cells.Item(19,6).Value = "NO"

# Row 20: match 4673
$battingExtra.Cells.Item(20,1).NumberFormat = "@"
$battingExtra.Cells.Item(20,1).Value = "4673"
$battingExtra.Cells.Item(20,2).Value = ""
$battingExtra.Cells.Item(20,3).Value = ""
$battingExtra.Cells.Item(20,4).Value = ""
$battingExtra.Cells.Item(20,5).Value = ""
$battingExtra.Cells.Item(20,6).Value = "NO"

# Row 21: match 4676
$battingExtra.Cells.Item(21,1).NumberFormat = "@"
$battingExtra.Cells.Item(21,1).Value = "4676"
$battingExtra.Cells.Item(21,2).Value = 4
$battingExtra.Cells.Item(21,3).NumberFormat = "@"
$battingExtra.Cells.Item(21,3).Value = "2"
$battingExtra.Cells.Item(21,4).NumberFormat = "@"
$battingExtra.Cells.Item(21,4).Value = "0"
$battingExtra.Cells.Item(21,5).Value = "4.57%"
$battingExtra.Cells.Item(21,6).Value = "NO"

